$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Element Types" sheet -- rename two categories, add a header, and
#    re-sort the whole list alphabetically.
# ---------------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Element Types")

# "Class/Interface" -> "Classes/Interface"
$wsElem.Range("A3").Value = "Classes/Interface"

# New header for the helper column.
$wsElem.Range("A1").Value = "Element Types"

# Re-sort the remaining (already-existing) category labels alphabetically.
$wsElem.Range("A4").Value = "Classic BAdI Implementation"
$wsElem.Range("A5").Value = "Custom Fiori Application"
$wsElem.Range("A6").Value = "Enhancement Implementation"
$wsElem.Range("A7").Value = "Extra Workbench Object"
$wsElem.Range("A8").Value = "Function Group"
$wsElem.Range("A9").Value = "Function Module"
$wsElem.Range("A10").Value = "Package"
$wsElem.Range("A11").Value = "Program"
$wsElem.Range("A12").Value = "Table"
$wsElem.Range("A13").Value = "Transaction"

# "BTP App / Extension" -> "BTP App/Extension" (written last so this is the
# final brand-new shared string interned).
$wsElem.Range("A2").Value = "BTP App/Extension"

# ---------------------------------------------------------------------------
# 2. "Library Elements" sheet
# ---------------------------------------------------------------------------
$wsLib = $wb.Worksheets.Item("Library Elements")

# Remove the stale review comment on E3 ("System Group, Component not
# allowed for configuration activities - ignore or warn").
$commentE3 = $wsLib.Range("E3").Comment
if ($commentE3) {
    $commentE3.Delete()
}

# Update the sample "Element Type" values in the demo rows.
$wsLib.Range("G3").Value = "BTP App/Extension"
$wsLib.Range("G4").Value = "Program"

# Add a dropdown (list data validation) on the Element Type column driven by
# the "Element Types" helper sheet.
$wsLib.Range("G2:G10").Validation.Add(3, 1, 1, "='Element Types'!`$A`$2:`$A`$13")

# ---------------------------------------------------------------------------
# 3. "Library Type" sheet -- drop the empty, unused third column.
# ---------------------------------------------------------------------------
$wsType = $wb.Worksheets.Item("Library Type")
$wsType.Columns.Item(2).Delete()
